# Swap the presentation's applied theme colour scheme ("Integral") for the
# classic default Office theme colours.
#
# The deck's slide master / overall design resolve to a single shared theme
# part; touching any slide's ThemeColorScheme mutates that one shared part,
# so we only need to set it once (accessed via slide 1) to recolour the
# whole deck.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$tcs.Item(1).RGB  = 0x000000   # dk1
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1
$tcs.Item(3).RGB  = 0x6A5444   # dk2      (44546A)
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      (E7E6E6)
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  (5B9BD5)
$tcs.Item(6).RGB  = 0x317DED   # accent2  (ED7D31)
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  (A5A5A5)
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  (FFC000)
$tcs.Item(9).RGB  = 0xC47244   # accent5  (4472C4)
$tcs.Item(10).RGB = 0x47AD70   # accent6  (70AD47)
$tcs.Item(11).RGB = 0xC16305   # hlink    (0563C1)
$tcs.Item(12).RGB = 0x724F95   # folHlink (954F72)
